$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# New rows 136-148 data: A (string), B, C, D (=IF(B>0,C-B,0))
$rows = @(
    @(136, "1st Move",        45752, 53228),
    @(137, "Checkpoint 89",   45830, 53306),
    @(138, "Checkpoint 404",  45929, 53405),
    @(139, "Checkpoint 1037", 46124, 53602),
    @(140, "Checkpoint 1534", 46279, 53758),
    @(141, "Checkpoint 1836", 46376, 53854),
    @(142, "Checkpoitn 2224", 46504, 53982),
    @(143, "Checkpoint 2586", 46624, 54102),
    @(144, "Enter door",      46876, 54354),
    @(145, "Touch button",    47892, 55387),
    @(146, "End level",       48805, 56300),
    @(147, "Enter 8-5",       50431, 59573),
    @(148, "1st Move",        50658, 59821)
)

foreach ($row in $rows) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $c = $row[3]

    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Formula = "=IF(B$r > 0,C$r-B$r, 0)"
}

# Update the frozen pane / selection view to match the new sheet extent.
# (topLeftCell of the frozen pane tracks the freeze boundary in this host;
# what we can control directly is the final selected cell.)
$ws.Activate()
$ws.Range("B149").Select()
